$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking Price cells as Text so they keep their exact
# string representation (e.g. "1.00", "130.00") instead of Excel auto-
# converting them to numbers (1, 130) as it would for a plain numeric entry.
$textCells = @("D4","D5","D6","D9","D10","D11","D12","D13","D15","D18","D19","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D33","D34","D36","D37","D42","D43","D45","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}


# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.074.51"
$ws.Range("E2").Value = "  -0.39%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.431.00"
$ws.Range("E3").Value = "  -0.01%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "409.39"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6 - Solana
$ws.Range("D6").Value = "130.00"
$ws.Range("E6").Value = "  -2.24%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +6.14%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.08%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.743"
$ws.Range("E9").Value = "  +7.21%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  +6.27%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "42.94"

# Row 12 - ShibaInu
$ws.Range("D12").Value = "0.0000230"
$ws.Range("E12").Value = "  +54.58%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "9.25"
$ws.Range("E13").Value = "  +10.04%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.17%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "21.45"
$ws.Range("E15").Value = "  +7.88%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "3.973.36"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.453.34"
$ws.Range("E17").Value = "  +0.78%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "12.51"
$ws.Range("E18").Value = "  +6.85%  "

# Row 19 - Polygon
$ws.Range("D19").Value = "1.10"
$ws.Range("E19").Value = "  +7.89%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "62.116.96"
$ws.Range("E20").Value = "  -0.27%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "453.53"
$ws.Range("E21").Value = "  +45.47%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "91.71"
$ws.Range("E22").Value = "  +8.97%  "

# Row 23 - ImmutableX
$ws.Range("D23").Value = "3.22"
$ws.Range("E23").Value = "  +1.66%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "13.05"
$ws.Range("E24").Value = "  +2.36%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +2.42%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "33.07"
$ws.Range("E26").Value = "  +11.21%  "

# Row 27 - Filecoin
$ws.Range("D27").Value = "8.98"
$ws.Range("E27").Value = "  +10.00%  "

# Row 28 - LEO
$ws.Range("D28").Value = "4.79"
$ws.Range("E28").Value = "  +1.22%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "7.74"
$ws.Range("E29").Value = "  -0.73%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "2.70"
$ws.Range("E30").Value = "  -1.95%  "

# Row 31 - Cosmos
$ws.Range("D31").Value = "12.02"
$ws.Range("E31").Value = "  +5.93%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  -1.11%  "

# Row 33 - InjectiveProtocol -> Hedera (rows 33/34 swapped)
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  -0.32%  "

# Row 34 - Hedera -> InjectiveProtocol (rows 33/34 swapped)
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "42.97"
$ws.Range("E34").Value = "  -1.45%  "

# Row 35 - Dai
$ws.Range("E35").Value = "  -0.09%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "0.0503"
$ws.Range("E36").Value = "  +3.45%  "

# Row 37 - OKB
$ws.Range("D37").Value = "54.29"
$ws.Range("E37").Value = "  +5.14%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.10%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +1.52%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +7.39%  "

# Row 41 - TheGraph
$ws.Range("E41").Value = "  +1.66%  "

# Row 42 - Stacks
$ws.Range("D42").Value = "2.95"
$ws.Range("E42").Value = "  -1.75%  "

# Row 43 - Monero
$ws.Range("D43").Value = "143.01"
$ws.Range("E43").Value = "  +0.28%  "

# Row 44 - NEARProtocol
$ws.Range("E44").Value = "  +8.98%  "

# Row 45 - ARBITRUM
$ws.Range("D45").Value = "2.00"
$ws.Range("E45").Value = "  +0.97%  "

# Row 46 - WEMIXToken
$ws.Range("E46").Value = "  +13.31%  "

# Row 47 - Celestia
$ws.Range("D47").Value = "16.71"
$ws.Range("E47").Value = "  -0.68%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "22.47"
$ws.Range("E48").Value = "  +5.04%  "

# Row 49 - ThetaToken
$ws.Range("D49").Value = "2.12"
$ws.Range("E49").Value = "  +8.75%  "

# Row 50 - RocketPoolETH -> Cronos (rows 50/51 swapped)
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.139"
$ws.Range("E50").Value = "  +16.19%  "

# Row 51 - Cronos -> RocketPoolETH (rows 50/51 swapped)
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "3.779.28"
$ws.Range("E51").Value = "  -0.22%  "
